{"js": "// Populate the intra-assay and inter-assay variability tables with Jinja\n// template placeholders (inter_var_sampleN_*) and drop the per-cell center\n// justification that Word had been writing on every paragraph in those\n// tables' data rows (sample-number, n, mean, SD, CV).\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Within this template, table index 4 is the \"intra-assay precision\" table\n// and table index 5 is the \"inter-assay precision\" table. Both share the\n// exact same 5-column layout (Sample | n | Mean | SD | CV) and, per the\n// target edit, both are populated with the *same* set of template variable\n// names (inter_var_sampleN_n / _mean / _sd / _cv).\nconst targetTables = [tables.items[4], tables.items[5]];\n\n// Load all rows/cells up front.\nfor (const tbl of targetTables) {\n  tbl.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const tbl of targetTables) {\n  for (const row of tbl.rows.items) {\n    row.cells.load(\"items\");\n  }\n}\nawait context.sync();\n\n// Column order after the \"Sample\" column: n, Mean, Standard Deviation, CV.\nconst suffixes = [\"n\", \"mean\", \"sd\", \"cv\"];\n\nfor (const tbl of targetTables) {\n  const rows = tbl.rows.items;\n  // Row 0 is the bold header row (Sample/n/Mean/SD/CV) and is left as-is.\n  for (let ri = 1; ri < rows.length; ri++) {\n    const cells = rows[ri].cells.items;\n    const sampleNum = ri; // data rows are Sample 1, 2, 3 in order\n\n    // Column 0: the \"Sample\" number cell keeps its literal text, but still\n    // loses the centered paragraph formatting.\n    cells[0].body.load(\"text\");\n  }\n}\nawait context.sync();\n\nfor (const tbl of targetTables) {\n  const rows = tbl.rows.items;\n  for (let ri = 1; ri < rows.length; ri++) {\n    const cells = rows[ri].cells.items;\n    const sampleNum = ri;\n\n    // Re-insert the cell's own text to flush out the centered <w:pPr> while\n    // keeping the displayed value identical.\n    cells[0].body.insertText(cells[0].body.text, \"Replace\");\n\n    // Columns 1..4: n / mean / sd / cv -> template placeholders.\n    for (let ci = 1; ci < cells.length; ci++) {\n      const varName = `inter_var_sample${sampleNum}_${suffixes[ci - 1]}`;\n      cells[ci].body.insertText(`{{ ${varName} }}`, \"Replace\");\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Populate the intra-assay and inter-assay variability tables with Jinja\n# template placeholders (inter_var_sampleN_*) and drop the per-cell center\n# justification that had been applied to every paragraph in those tables'\n# data rows (sample-number, n, mean, SD, CV).\n\n$d = $word.ActiveDocument\n\n# Table index 4 (1-based -> Tables.Item(5)) is the \"intra-assay precision\"\n# table and table index 5 (1-based -> Tables.Item(6)) is the \"inter-assay\n# precision\" table. Both share the exact same 5-column layout\n# (Sample | n | Mean | SD | CV) and, per the target edit, both get populated\n# with the *same* set of template variable names\n# (inter_var_sampleN_n / _mean / _sd / _cv).\n$tableIndexes = @(5, 6)\n\n# Column order after the \"Sample\" column: n, Mean, Standard Deviation, CV.\n$suffixes = @(\"n\", \"mean\", \"sd\", \"cv\")\n\nforeach ($tableIndex in $tableIndexes) {\n    $t = $d.Tables.Item($tableIndex)\n\n    # Row 1 is the bold header row (Sample/n/Mean/SD/CV) and is left as-is.\n    for ($rowIndex = 2; $rowIndex -le $t.Rows.Count; $rowIndex++) {\n        $sampleNum = $rowIndex - 1\n\n        # Column 1: the \"Sample\" number cell keeps its literal text, but\n        # still loses the centered paragraph formatting.\n        $sampleCell = $t.Cell($rowIndex, 1)\n        $sampleCell.Range.ParagraphFormat.Alignment = 0\n\n        # Columns 2..5: n / mean / sd / cv -> template placeholders.\n        for ($colIndex = 2; $colIndex -le $t.Columns.Count; $colIndex++) {\n            $varName = \"inter_var_sample\" + $sampleNum + \"_\" + $suffixes[$colIndex - 2]\n            $cell = $t.Cell($rowIndex, $colIndex)\n            $cell.Range.Text = \"{{ \" + $varName + \" }}\"\n            $cell.Range.ParagraphFormat.Alignment = 0\n        }\n    }\n}\n"}
